$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# The repayment schedule gained an extra (currently blank) column between the
# existing "In Advance" (M) and "Late" (N) columns, to support the new
# "Variable Instalments" data point. Inserting a column shifts "Late" from N
# to O and "Outstanding" from O/P to P/Q, leaving the new N column blank -
# exactly matching the target layout.
$ws.Columns("N").Insert() | Out-Null

# Re-select the cell that was active when the workbook was last saved.
$ws.Range("S5").Select() | Out-Null
